$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.793.46'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.353.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.54%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.74%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.61'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.600'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.61%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.36'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '37.22'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +16.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.37'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.60%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.702.89'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.59%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.926'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.346.99'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.737.81'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.05%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.16'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.15'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.94%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.87'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.96%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.71'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.30'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.93'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.134'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.56'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.26%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.00%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.39'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.95%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.71%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.33'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +14.60%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +18.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.25'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +9.66%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.98%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.203'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.41%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.08'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.09%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.70%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.24'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.28%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.19'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.00%  '
